# B6-PowerPoint.pptx edit:
#  1) Re-apply the table style on the three summary tables (slides 14-16)
#     so they point at the built-in style {C78438C0-87C2-462F-93DF-041E0640FE46}
#     instead of the custom "Table_0" style that used to be referenced.
#  2) The deck ships two theme parts: theme2.xml (linked from the slide
#     master - the "Integral"/Red Violet palette actually seen on the
#     slides) and theme1.xml (linked only from the notes master - a
#     plain "Office Theme" palette). The edit swaps their two colour
#     palettes. The slide-master palette (theme2.xml) is reachable from
#     the object model via SlideMaster.Theme, so push the "Office Theme"
#     colours into it, in clrScheme order: dk1, lt1, dk2, lt2,
#     accent1-6, hlink, folHlink. (COM RGB longs are 0xBBGGRR.)

$p = $ppt.ActivePresentation

$newStyle = "{C78438C0-87C2-462F-93DF-041E0640FE46}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyle)
    }
}

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0x000000   # dk1      000000
$colorScheme.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colorScheme.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colorScheme.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colorScheme.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colorScheme.Item(11).RGB = 0xC16305   # hlink    0563C1
$colorScheme.Item(12).RGB = 0x724F95   # folHlink 954F72
